# NIT-9008957041 Estado de Cuenta update:
#  - Previous account-statement record for worker ALEJANDRA GARCIA ALIAN (row 17) is removed.
#  - The "Periodo Mora" for the remaining worker (LUIS ANGEL SAYAS RUIZ) moves from 2506 to 2507.
#  - Totals are refreshed to reflect the single remaining worker: Valor Mora 113880 -> 56940,
#    Cant. Trabajadores 2 -> 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for the second worker (ALEJANDRA GARCIA ALIAN), shifting rows below up.
$ws.Rows.Item(17).Delete()

# Update the remaining worker's "Periodo Mora" (kept as text, matching original cell type).
$ws.Range("E16").Value = "2507"

# Refresh the summary totals now that only one worker/record remains.
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
